$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Approved/Rejected (I) + ReasonToReject (J) go from
# "Rejected"/"Nil" -> "Approved" only (J cell cleared).
$toApproved = @(2, 4, 5, 7, 8, 10, 12, 33, 49, 56)
foreach ($r in $toApproved) {
    $ws.Range("J$r").ClearContents()
    $ws.Range("I$r").Value2 = "Approved"
}

# Rows whose Approved/Rejected (I) go from "Approved" only ->
# "Rejected"/"Nil", highlighted with a yellow fill.
$toRejected = @(40, 47, 58, 65, 67, 75, 77, 85, 87, 95, 97)
foreach ($r in $toRejected) {
    $ws.Range("I$r").Value2 = "Rejected"
    $ws.Range("J$r").Value2 = "Nil"
    $ws.Range("I$r`:J$r").Interior.Color = 65535
}

# Move the active cell selection from I56 to J56.
$ws.Range("J56").Select()
